$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 137929
$ws.Range("D3").Value = 140306
$ws.Range("D4").Value = 143475
$ws.Range("D5").Value = 145801
$ws.Range("D6").Value = 147852
$ws.Range("D7").Value = 153419
$ws.Range("D8").Value = 154766
$ws.Range("D9").Value = 145498
